# Updated cryptos list values (Price / Volume(1h), plus the
# HuobiToken/NEARProtocol row swap) per the target commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.495.99'
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").Value = '2.160.65'
$ws.Range("E3").Value = '  -2.90%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'238.30"
$ws.Range("E5").Value = '  -1.98%  '
$ws.Range("D6").Value = "'0.606"
$ws.Range("E6").Value = '  -3.28%  '
$ws.Range("D7").Value = "'71.80"
$ws.Range("E7").Value = '  -2.17%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.576"
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").Value = "'39.76"
$ws.Range("E10").Value = '  -6.20%  '
$ws.Range("D11").Value = "'0.0904"
$ws.Range("E11").Value = '  -5.30%  '
$ws.Range("D12").Value = "'54.13"
$ws.Range("E12").Value = '  -4.58%  '
$ws.Range("D13").Value = "'0.0999"
$ws.Range("E13").Value = '  -3.82%  '
$ws.Range("D14").Value = "'6.69"
$ws.Range("E14").Value = '  -3.47%  '
$ws.Range("D15").Value = '2.484.47'
$ws.Range("E15").Value = '  -2.98%  '
$ws.Range("D16").Value = "'14.05"
$ws.Range("E16").Value = '  -1.93%  '
$ws.Range("D17").Value = '2.157.74'
$ws.Range("E17").Value = '  -3.96%  '
$ws.Range("D18").Value = "'0.779"
$ws.Range("E18").Value = '  -7.07%  '
$ws.Range("D19").Value = '41.347.93'
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("E20").Value = '  -2.93%  '
$ws.Range("D21").Value = "'69.58"
$ws.Range("E21").Value = '  -4.47%  '
$ws.Range("D22").Value = "'5.76"
$ws.Range("E22").Value = '  -7.91%  '
$ws.Range("D23").Value = "'9.82"
$ws.Range("E23").Value = '  -12.89%  '
$ws.Range("D24").Value = "'227.14"
$ws.Range("D25").Value = "'2.01"
$ws.Range("E25").Value = '  -3.50%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = "'10.64"
$ws.Range("E27").Value = '  -6.56%  '
$ws.Range("D28").Value = "'3.29"
$ws.Range("E28").Value = '  -9.29%  '
$ws.Range("D29").Value = "'2.18"
$ws.Range("E29").Value = '  -4.59%  '
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").Value = "'169.74"
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("D32").Value = "'19.74"
$ws.Range("E32").Value = '  -3.81%  '
$ws.Range("D33").Value = "'33.29"
$ws.Range("E33").Value = '  +10.68%  '
$ws.Range("E34").Value = '  -4.18%  '
$ws.Range("D35").Value = "'5.14"
$ws.Range("E35").Value = '  -8.19%  '
$ws.Range("E36").Value = '  -3.96%  '
$ws.Range("D37").Value = "'0.104"
$ws.Range("E37").Value = '  -3.67%  '
$ws.Range("D38").Value = "'4.28"
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").Value = "'0.0299"
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").Value = "'12.09"
$ws.Range("E40").Value = '  -10.17%  '
$ws.Range("D41").Value = "'2.07"
$ws.Range("E41").Value = '  -2.92%  '
$ws.Range("D42").Value = "'5.31"
$ws.Range("E42").Value = '  -6.10%  '
$ws.Range("D43").Value = "'58.85"
$ws.Range("E43").Value = '  -9.34%  '
$ws.Range("D44").Value = "'0.188"
$ws.Range("E44").Value = '  -4.64%  '
$ws.Range("D45").Value = "'8.34"
$ws.Range("E45").Value = '  -4.65%  '
$ws.Range("D46").Value = "'0.0957"
$ws.Range("E46").Value = '  -4.50%  '
$ws.Range("D47").Value = "'95.85"
$ws.Range("E47").Value = '  -8.80%  '
$ws.Range("E48").Value = '  -3.33%  '
$ws.Range("D49").Value = "'1.11"
$ws.Range("E49").Value = '  -5.21%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = "'2.17"
$ws.Range("E50").Value = '  -7.62%  '
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").Value = "'2.62"
$ws.Range("E51").Value = '  -2.68%  '
